$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.200.27"
$ws.Range("D3").Value = "1.858.10"
$ws.Range("E3").Value = "  +1.69%  "
$ws.Range("E4").Value = "  +0.51%  "
$ws.Range("D5").Value = "'239.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'0.622"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.75%  "
$ws.Range("E7").Value = "  +0.49%  "
$ws.Range("D8").Value = "'42.22"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +6.95%  "
$ws.Range("E9").Value = "  +0.99%  "
$ws.Range("E10").Value = "  +1.59%  "
$ws.Range("D11").Value = "'0.0990"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.48%  "
$ws.Range("D12").Value = "2.127.53"
$ws.Range("E12").Value = "  +1.70%  "
$ws.Range("D13").Value = "'11.50"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.87%  "
$ws.Range("D14").Value = "1.864.87"
$ws.Range("E14").Value = "  +1.79%  "
$ws.Range("D15").Value = "'0.676"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.08%  "
$ws.Range("E16").Value = "  +2.37%  "
$ws.Range("D17").Value = "35.197.07"
$ws.Range("E17").Value = "  +1.25%  "
$ws.Range("E18").Value = "  +0.63%  "
$ws.Range("E19").Value = "  +1.42%  "
$ws.Range("D20").Value = "'240.77"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.19%  "
$ws.Range("E21").Value = "  +0.94%  "
$ws.Range("D22").Value = "'4.76"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.93%  "
$ws.Range("E23").Value = "  +0.45%  "
$ws.Range("D24").Value = "'2.27"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.07%  "
$ws.Range("D25").Value = "'169.77"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.95%  "
$ws.Range("E26").Value = "  +26.42%  "
$ws.Range("E27").Value = "  +3.49%  "
$ws.Range("D28").Value = "'17.67"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.04%  "
$ws.Range("E29").Value = "  +0.24%  "
$ws.Range("D30").Value = "'0.0561"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.21%  "
$ws.Range("E31").Value = "  +0.49%  "
$ws.Range("E32").Value = "  +2.19%  "
$ws.Range("D33").Value = "'1.82"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +26.79%  "
$ws.Range("E34").Value = "  +2.23%  "
$ws.Range("E35").Value = "  +10.76%  "
$ws.Range("D36").Value = "'0.820"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +17.89%  "
$ws.Range("E37").Value = "  +7.88%  "
$ws.Range("D38").Value = "'1.10"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.31%  "
$ws.Range("E39").Value = "  +4.58%  "
$ws.Range("D40").Value = "'89.92"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.16%  "
$ws.Range("D41").Value = "1.347.48"
$ws.Range("D42").Value = "'0.0603"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +15.61%  "
$ws.Range("E43").Value = "  +3.32%  "
$ws.Range("E44").Value = "  +3.13%  "
$ws.Range("E45").Value = "  +0.45%  "
$ws.Range("D46").Value = "'12.52"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +42.60%  "
$ws.Range("E47").Value = "  -0.38%  "
$ws.Range("D48").Value = "'6.57"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.36%  "
$ws.Range("D49").Value = "2.046.01"
$ws.Range("E49").Value = "  +1.94%  "
$ws.Range("E50").Value = "  +1.37%  "
$ws.Range("E51").Value = "  +0.56%  "
